$d = $word.ActiveDocument

# 1. Merge the "Facebook" user story runs into a single run.
$d.Content.Find.Execute(
    "As an unauthorized user I want to use my Facebook account so that I can use one of my own existing accounts.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As an unauthorized user I want to use my Facebook account so that I can use one of my own existing accounts.",
    2)

# 2. Merge the "CSwap" user story runs (and drop the proofErr spell-check markers) into a single run.
$d.Content.Find.Execute(
    "As a student looking for furniture, I want to use CSwap and navigate to the furniture section so that I can find furniture for my apartment.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As a student looking for furniture, I want to use CSwap and navigate to the furniture section so that I can find furniture for my apartment.",
    2)

# 3. Split the "poor user" run into three runs, replacing "poor" with "student".
$find = $d.Content.Find
$find.Execute("poor")
$r = $find.Parent
$r.Font.Bold = $true
$r.Text = "student"
$r.Font.Bold = $false
